$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Abubakar Fatima Shehu"
$ws.Range("B3").Value = "std356"
$ws.Range("C3").Value = "SS3_GOLD"
$ws.Range("D3").Value = "CHEMISTRY"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2%"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = "FAIL"
$ws.Range("I3").Value = "2025-12-04 07:10"
